# Data updated by GitHub Bot (2020-06-03 12)
# Appends 9 new daily rows (43976..43984 -> 25 May 2020 .. 2 Jun 2020) to the
# "Covid-19 podatki" sheet, extends the Tabela1 table / AutoFilter to the new
# range, and moves the active selection to the last appended row - matching
# the upstream data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 75
$firstNewRow = 76
$lastNewRow = 84

# New daily figures: Date(serial), Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, Intensive care, Discharged, Deaths(all),
# Deaths(daily) - same column order as the worksheet / table.
$data = @(
    @(43976,75770,754,1469,0,9,2,6,108,1),
    @(43977,76579,809,1471,2,8,2,2,108,0),
    @(43978,77210,631,1473,2,7,2,1,108,0),
    @(43979,77916,706,1473,0,7,2,0,108,0),
    @(43980,78529,613,1473,0,7,2,0,108,0),
    @(43981,78793,264,1473,0,6,2,1,108,0),
    @(43982,79039,246,1473,0,5,1,0,109,1),
    @(43983,79698,659,1475,2,5,1,0,109,0),
    @(43984,80505,807,1477,2,5,0,0,109,0)
)

# Carry the formatting (number format / font / alignment) of the last data
# row down onto the new rows before writing values, otherwise new cells pick
# up the bare column default style (plain text) instead of the table's
# numeric/date formatting.
$fmtSource = $ws.Range("A" + $lastExistingRow + ":J" + $lastExistingRow)
$fmtTarget = $ws.Range("A" + $firstNewRow + ":J" + $lastNewRow)
$fmtSource.Copy($fmtTarget)

$r = $firstNewRow
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# Grow the table (and its AutoFilter) so it covers the freshly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J" + $lastNewRow))

# Match the author's final on-screen selection (last appended row).
$null = $ws.Range("A" + $lastNewRow + ":J" + $lastNewRow).Select()
